# Weekly update: insert the newest "Locoto" price-report rows (Primera /
# Segunda quality) right above the existing table body, pushing the
# historical rows down by two (the table keeps growing every week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 40 (old rows 40..72 shift to 42..74).
$ws.Range("A40:A41").EntireRow.Insert()

# --- New row 40: Primera ---
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44484
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112042
$ws.Range("G40").Value = "Locoto"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 140
$ws.Range("K40").Value = 26000
$ws.Range("L40").Value = 27000
$ws.Range("M40").Value = 26500
$ws.Range("N40").Value = "$/caja 20 kilos"
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 1325
$ws.Range("Q40").Value = 20
$ws.Range("R40").Value = "Hortaliza"

# --- New row 41: Segunda ---
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 44484
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112042
$ws.Range("G41").Value = "Locoto"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Segunda"
$ws.Range("J41").Value = 120
$ws.Range("K41").Value = 24000
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = 24500
$ws.Range("N41").Value = "$/caja 20 kilos"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 1225
$ws.Range("Q41").Value = 20
$ws.Range("R41").Value = "Hortaliza"
